$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (Ano 2025) with refreshed faturamento data
$ws.Range("B9").Value = 3588652.36
$ws.Range("C9").Value = 564405.24
$ws.Range("D9").Value = 4153057.6
$ws.Range("E9").Value = 13.59011346242826
$ws.Range("F9").Value = 86.40988653757175
$ws.Range("G9").Value = -45.45304780826843
$ws.Range("H9").Value = -35.19389875585478
$ws.Range("I9").Value = 36170
$ws.Range("J9").Value = 1546
$ws.Range("K9").Value = 37716
$ws.Range("L9").Value = 26050
$ws.Range("M9").Value = 159.4263953934741
$ws.Range("N9").Value = 8.84364728880691
